# Add a "Remarks" column (K) to the shareholder upload template and
# adjust the frozen panes so the first two columns (A:B) stay pinned
# alongside the header row when scrolling, matching the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Remarks" header in column K -------------------------------------
$ws.Range("K1").Value = "Remarks"

# Match the formatting of the existing header cells (border + text format)
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Widen the new column like the other header columns
$ws.Columns.Item(11).ColumnWidth = 27.59

# --- Freeze panes at C2 (freeze header row + first two columns) -----------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("C2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("A1").Select() | Out-Null
